# Publish terminology IG 2.0.2
# Updates the Metadata sheet: Version, Status, Experimental (cleared), Date

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.8.1 -> 1.8.2
$ws.Range("B3").Value = "1.8.2"

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: true -> (cleared)
$ws.Range("B7").ClearContents()

# Date: 2023-10-31 -> 2025-11-18 (keep as text, not a date serial)
$ws.Range("B8").Formula = "=""2025-11-18"""
$ws.Range("B8").Copy()
$ws.Range("B8").PasteSpecial(-4163)
